$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-53 with refreshed date (A) and simulated value (B) data
$ws.Cells.Item(2, 1).Value = [double]"39400"
$ws.Cells.Item(2, 2).Value = [double]"-0.03834288659695062"
$ws.Cells.Item(3, 1).Value = [double]"39583"
$ws.Cells.Item(3, 2).Value = [double]"1.257433230729447"
$ws.Cells.Item(4, 1).Value = [double]"39765"
$ws.Cells.Item(4, 2).Value = [double]"0.8837904892317567"
$ws.Cells.Item(5, 1).Value = [double]"39948"
$ws.Cells.Item(5, 2).Value = [double]"0.2524978494830066"
$ws.Cells.Item(6, 1).Value = [double]"40130"
$ws.Cells.Item(6, 2).Value = [double]"0.5778148852415939"
$ws.Cells.Item(7, 1).Value = [double]"40310"
$ws.Cells.Item(7, 2).Value = [double]"-0.1977656654399595"
$ws.Cells.Item(8, 1).Value = [double]"40494"
$ws.Cells.Item(8, 2).Value = [double]"-1.091011900795806"
$ws.Cells.Item(9, 1).Value = [double]"40676"
$ws.Cells.Item(9, 2).Value = [double]"-0.1775928823643795"
$ws.Cells.Item(10, 1).Value = [double]"40862"
$ws.Cells.Item(10, 2).Value = [double]"0.5531759638372762"
$ws.Cells.Item(11, 1).Value = [double]"41044"
$ws.Cells.Item(11, 2).Value = [double]"0.1947850960503388"
$ws.Cells.Item(12, 1).Value = [double]"41228"
$ws.Cells.Item(12, 2).Value = [double]"0.3732050716642448"
$ws.Cells.Item(13, 1).Value = [double]"41409"
$ws.Cells.Item(13, 2).Value = [double]"-0.1411005862636046"
$ws.Cells.Item(14, 1).Value = [double]"41592"
$ws.Cells.Item(14, 2).Value = [double]"-0.2128461555332564"
$ws.Cells.Item(15, 1).Value = [double]"41774"
$ws.Cells.Item(15, 2).Value = [double]"0.4255979180752121"
$ws.Cells.Item(16, 1).Value = [double]"41957"
$ws.Cells.Item(16, 2).Value = [double]"0.4461687925667093"
$ws.Cells.Item(17, 1).Value = [double]"42137"
$ws.Cells.Item(17, 2).Value = [double]"0.3261422475203943"
$ws.Cells.Item(18, 1).Value = [double]"42321"
$ws.Cells.Item(18, 2).Value = [double]"0.6601374471387373"
$ws.Cells.Item(19, 1).Value = [double]"42503"
$ws.Cells.Item(19, 2).Value = [double]"0.8728685839363095"
$ws.Cells.Item(20, 1).Value = [double]"42689"
$ws.Cells.Item(20, 2).Value = [double]"0.9596379771730028"
$ws.Cells.Item(21, 1).Value = [double]"42867"
$ws.Cells.Item(21, 2).Value = [double]"0.2682953781150843"
$ws.Cells.Item(22, 1).Value = [double]"43053"
$ws.Cells.Item(22, 2).Value = [double]"0.1682050168937224"
$ws.Cells.Item(23, 1).Value = [double]"43145"
$ws.Cells.Item(23, 2).Value = [double]"0.4761878885828992"
$ws.Cells.Item(24, 1).Value = [double]"43235"
$ws.Cells.Item(24, 2).Value = [double]"-0.5240674734835977"
$ws.Cells.Item(25, 1).Value = [double]"43326"
$ws.Cells.Item(25, 2).Value = [double]"-0.2614280992174685"
$ws.Cells.Item(26, 1).Value = [double]"43418"
$ws.Cells.Item(26, 2).Value = [double]"0.8"
$ws.Cells.Item(27, 1).Value = [double]"43510"
$ws.Cells.Item(27, 2).Value = [double]"-0.2951285663450562"
$ws.Cells.Item(28, 1).Value = [double]"43600"
$ws.Cells.Item(28, 2).Value = [double]"-0.3"
$ws.Cells.Item(29, 1).Value = [double]"43691"
$ws.Cells.Item(29, 2).Value = [double]"0.7916058519991367"
$ws.Cells.Item(30, 1).Value = [double]"43783"
$ws.Cells.Item(30, 2).Value = [double]"0.5670926739443871"
$ws.Cells.Item(31, 1).Value = [double]"43875"
$ws.Cells.Item(31, 2).Value = [double]"1.261225886527512"
$ws.Cells.Item(32, 1).Value = [double]"43966"
$ws.Cells.Item(32, 2).Value = [double]"0.2"
$ws.Cells.Item(33, 1).Value = [double]"44068"
$ws.Cells.Item(33, 2).Value = [double]"1.450185044412038"
$ws.Cells.Item(34, 1).Value = [double]"44159"
$ws.Cells.Item(34, 2).Value = [double]"0.760290197271857"
$ws.Cells.Item(35, 1).Value = [double]"44251"
$ws.Cells.Item(35, 2).Value = [double]"-0.5062757877985717"
$ws.Cells.Item(36, 1).Value = [double]"44341"
$ws.Cells.Item(36, 2).Value = [double]"0.1999914549530217"
$ws.Cells.Item(37, 1).Value = [double]"44432"
$ws.Cells.Item(37, 2).Value = [double]"1.809651846369789"
$ws.Cells.Item(38, 1).Value = [double]"44525"
$ws.Cells.Item(38, 2).Value = [double]"-2.218509177329452"
$ws.Cells.Item(39, 1).Value = [double]"44617"
$ws.Cells.Item(39, 2).Value = [double]"0.9580856546984649"
$ws.Cells.Item(40, 1).Value = [double]"44706"
$ws.Cells.Item(40, 2).Value = [double]"0.06845131779844849"
$ws.Cells.Item(41, 1).Value = [double]"44798"
$ws.Cells.Item(41, 2).Value = [double]"2.344079640195559"
$ws.Cells.Item(42, 1).Value = [double]"44890"
$ws.Cells.Item(42, 2).Value = [double]"-5.587651344285405E-06"
$ws.Cells.Item(43, 1).Value = [double]"44981"
$ws.Cells.Item(43, 2).Value = [double]"0.6233623182505568"
$ws.Cells.Item(44, 1).Value = [double]"45071"
$ws.Cells.Item(44, 2).Value = [double]"-4.887355777631228"
$ws.Cells.Item(45, 1).Value = [double]"45163"
$ws.Cells.Item(45, 2).Value = [double]"0.07746478018819403"
$ws.Cells.Item(46, 1).Value = [double]"45254"
$ws.Cells.Item(46, 2).Value = [double]"0.2231996868496964"
$ws.Cells.Item(47, 1).Value = [double]"45345"
$ws.Cells.Item(47, 2).Value = [double]"0.3473129037311367"
$ws.Cells.Item(48, 1).Value = [double]"45436"
$ws.Cells.Item(48, 2).Value = [double]"-0.4268782796002455"
$ws.Cells.Item(49, 1).Value = [double]"45534"
$ws.Cells.Item(49, 2).Value = [double]"1.034675296103259"
$ws.Cells.Item(50, 1).Value = [double]"45618"
$ws.Cells.Item(50, 2).Value = [double]"0.434366938073353"
$ws.Cells.Item(51, 1).Value = [double]"45713"
$ws.Cells.Item(51, 2).Value = [double]"0.41600913674678"
$ws.Cells.Item(52, 1).Value = [double]"45800"
$ws.Cells.Item(52, 2).Value = [double]"-0.3423454266220887"
$ws.Cells.Item(53, 1).Value = [double]"45891"
$ws.Cells.Item(53, 2).Value = [double]"0.7919400257838731"

# Remove now-unused rows 54-73 (shrinks used range from B73 to B53)
$ws.Range("A54:B73").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

